# Adds team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF with values, formatted like the
# existing header cells (bold, centered, thin border - same style as AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-44 get the team's win/loss/tie record as numbers.
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 76   # AD
    $ws.Cells.Item($r, 31).Value = 86   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
